# Patch antenna workbook - add new revision (r0.1) data and tweak
# existing formulas to use absolute references so they still work
# when copied for the new revision block below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Relabel the first "Actual Design" section as revision r0.1
# ---------------------------------------------------------------
$ws.Range("A25").Value = "Actual Design r0.1"

# ---------------------------------------------------------------
# 2. Make the key formulas use absolute references (so that the
#    constants they pull from the top of the sheet keep resolving
#    correctly even though a second, near-identical block of
#    formulas is introduced lower on the sheet).
# ---------------------------------------------------------------
$ws.Range("B4").Formula  = "=(`$B`$3+1)/2+(`$B`$3-1)/2*((SQRT(1+12*`$B`$6/`$G`$14))^-1)"
$ws.Range("B8").Formula  = "=300000000/(`$B`$7*10^6)/SQRT(`$B`$4)"
$ws.Range("B9").Formula  = "=`$B`$8*0.0393701*1000"
$ws.Range("C9").Formula  = "=`$B`$9/2"
$ws.Range("D9").Formula  = "=`$B`$9/4"
$ws.Range("J13").Formula = "=G13/`$B`$9"
$ws.Range("J14").Formula = "=G14/`$B`$9"
$ws.Range("N27").Formula = "=G27/`$B`$9"
$ws.Range("N28").Formula = "=G28/`$B`$9"

# ---------------------------------------------------------------
# 3. Add the notes directly below the first "Actual Design" table
#    and introduce the new "Actual Design r0.2" revision table.
# ---------------------------------------------------------------
$ws.Range("B34").Font.Bold = $true
$ws.Range("B34").Value = "Resonance at about 2.2GHz and 90ohm impedance."

$ws.Range("A36").Value = "Actual Design r0.2"

# Header row for the new revision table (row 37), matching the
# style/labels of the row 26 header above.
$ws.Range("B37:O37").Font.Bold = $true
$ws.Range("B37").Value = "Ground (mm)"
$ws.Range("C37").Value = "Ground (mils)"
$ws.Range("D37").Value = "Substrate (mm)"
$ws.Range("E37").Value = "Substrate (mils)"
$ws.Range("F37").Value = "Patch (mm)"
$ws.Range("G37").Value = "Patch (mils)"
$ws.Range("H37").Value = "Feed (mm)"
$ws.Range("I37").Value = "Feed (mils)"
$ws.Range("J37").Value = "Feed Pos X (mm)"
$ws.Range("K37").Value = "Feed Pos X (mils)"
$ws.Range("L37").Value = "Feed Pos Y (mm)"
$ws.Range("M37").Value = "Feed Pos Y (mils)"
$ws.Range("O37").Value = "of effective wavelength"

# Row 38 - Width
$ws.Range("A38").Value = "Width"
$ws.Range("C38").Value = 1525
$ws.Range("B38").Formula = "=C38/0.0393701/1000"
$ws.Range("E38").Value = 1525
$ws.Range("D38").Formula = "=E38/0.0393701/1000"
$ws.Range("G38").Value = 1200
$ws.Range("F38").Formula = "=G38/0.0393701/1000"
$ws.Range("I38").Value = 1150
$ws.Range("H38").Formula = "=I38/0.0393701/1000"
$ws.Range("J38").Value = 1150
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("L38").Formula = "=M38/0.0393701/1000"
$ws.Range("N38").Font.Bold = $true
$ws.Range("N38").Formula = "=G38/`$B`$9"
$ws.Range("O38").Font.Bold = $true
$ws.Range("O38").Value = "of effective wavelength"

# Rows 39:40 - Length / Height, entered as ranges so that Excel
# stores them as shared formulas (si) just like the rest of the sheet.
$ws.Range("A39").Value = "Length"
$ws.Range("A40").Value = "Height"

$ws.Range("C39").Value = 1700
$ws.Range("C40").Value = 1.9
$ws.Range("B39:B40").Formula = "=C39/0.0393701/1000"

$ws.Range("E39").Value = 1700
$ws.Range("E40").Value = 57
$ws.Range("D39:D40").Formula = "=E39/0.0393701/1000"

$ws.Range("G39").Value = 1375
$ws.Range("G40").Value = 1.9
$ws.Range("F39:F40").Formula = "=G39/0.0393701/1000"

$ws.Range("I39").Value = 100
$ws.Range("I40").Value = 1.9
$ws.Range("H39:H40").Formula = "=I39/0.0393701/1000"

$ws.Range("K39").Value = 500
$ws.Range("K40").Value = 0
$ws.Range("J39:J40").Formula = "=K39/0.0393701/1000"

$ws.Range("M39").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("L39:L40").Formula = "=M39/0.0393701/1000"

$ws.Range("N39").Font.Bold = $true
$ws.Range("N39").Formula = "=G39/`$B`$9"

# ---------------------------------------------------------------
# 4. Notes for the new r0.2 revision.
# ---------------------------------------------------------------
$ws.Range("B42").Value = "Notes:"
$ws.Range("C42").Value = "Actual Design Notes r0.1"
$ws.Range("C43").Value = "Feed line should have the full width going into the signal net as opposed to a thin thermal relief."

# ---------------------------------------------------------------
# 5. Update the view so the new rows are visible / selected, mirroring
#    what Excel would have recorded after the edits were made.
# ---------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("A44").Select()
